$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New chronological data for rows 16-35 (Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico)
$data = @(
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1810", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1811", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1812", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1901", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1902", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1903", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1904", 31249, 781242),
    @("CC", "1082925993",  "GORTRUDE MARY MCLEAN CARDILES", "1905", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1905", 31249, 781242),
    @("CC", "1082925993",  "GORTRUDE MARY MCLEAN CARDILES", "1906", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1906", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1907", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1908", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1909", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1910", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1911", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "1912", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "2001", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "2002", 31249, 781242),
    @("CC", "73009373",    "XAVIER ELIAS TORRES MIRANDA",   "2003", 20833, 781242)
)

$startRow = 16
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
}
